# "send daily & weekly" — bump the reported week number (23 -> 24) for
# both data rows on the weekly template sheet, then leave the selection
# where the user's last edit landed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds the "week" value for each row (row 1 is the header).
$ws.Range("B2").Value = 24
$ws.Range("B3").Value = 24

# Move the active selection to B4, matching where the user's cursor
# ended up after editing the week column.
$ws.Range("B4").Select()
